$wb = $excel.ActiveWorkbook

# 1. Rename sheet "ResetEmployeeData" -> "ResetEmployeeData12"
$wsReset = $wb.Worksheets.Item("ResetEmployeeData")
$wsReset.Name = "ResetEmployeeData12"

# 2. Change the active sheet / selected tab to "ResetEmployeeData12"
#    and set its selection to D16.
$wsReset.Activate()
$wsReset.Range("D16").Select()
